# Update roomPerDay, roomArrangement, teacherArrangement
# Append a new "Fsoft" code/name row to the Department sheet and move the
# active selection to F10 (matching the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 11): Code = "Fsoft", Name = "Fsoft"
$ws.Range("A11").Value = "Fsoft"
$ws.Range("B11").Value = "Fsoft"

# Update the active selection to F10, as captured in the saved workbook
$ws.Range("F10").Select() | Out-Null
